# Adds a new "2022-Q3" quarter: updates the summary sheet ("总计") with a new
# leading row, and inserts a brand-new "2022-Q3" worksheet (placed right
# after "总计" and before "2022-Q2") containing the fund holdings detail for
# that quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value into a Range while forcing it to be stored as TEXT
# (keeps leading/trailing zeros exactly as typed, e.g. "000556" or "0.80"),
# then strip any leftover explicit formatting so the cell ends up with the
# default (no style) look, matching the rest of the plain data cells.
# ---------------------------------------------------------------------------
function Set-TextCell {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------------
# 1) Update the "总计" (summary) sheet: insert the 2022-Q3 row at the top of
#    the data (row 2) and shift the previously-existing rows down by one.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Existing rows 2..6 hold (in order): 2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3,
# 2021-Q1. They need to move to rows 3..7, so walk from the bottom up.
$existingRows = @(
    @{C = 6;  D = 0.31},                 # was row2: 2022-Q2
    @{C = 4;  D = 0.15},                  # was row3: 2022-Q1
    @{C = 7;  D = 0.26},                  # was row4: 2021-Q4
    @{C = 1;  D = 0.07000000000000001},   # was row5: 2021-Q3
    @{C = 18; D = 0.5}                    # was row6: 2021-Q1
)
$labels = @("2022-Q2", "2022-Q1", "2021-Q4", "2021-Q3", "2021-Q1")

# Give row 7's column A the same look (bold, centered, bordered) as the
# other index cells by copying the formatting from an existing one (A2).
$total.Range("A2").Copy() | Out-Null
$total.Range("A7").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

for ($i = 0; $i -lt 5; $i++) {
    $row = 7 - $i
    $total.Cells.Item($row, 1).Value = 5 - $i
    $total.Cells.Item($row, 2).Value = $labels[4 - $i]
    $total.Cells.Item($row, 3).Value = $existingRows[4 - $i].C
    $total.Cells.Item($row, 4).Value = $existingRows[4 - $i].D
}

# New first data row: 2022-Q3
$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 13
$total.Cells.Item(2, 4).Value = 0.21

# ---------------------------------------------------------------------------
# 2) Insert the brand-new "2022-Q3" worksheet right before "2022-Q2". Start
#    from a copy of "2022-Q2" so all sheet-level formatting/boilerplate
#    (margins, outline settings, header style, index-column style, ...)
#    matches the rest of the workbook, then wipe the copied values.
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"
$q3.Cells.ClearContents()

# The template only had 7 rows (6 funds); we need 14 (13 funds). Extend the
# bold/bordered "index" styling used in column A down to the extra rows.
$q3.Range("A2").Copy() | Out-Null
$q3.Range("A8:A14").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Header row (row 1), columns B..H - style already carried over by the copy.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $q3.Cells.Item(1, $col).Value = $headers[$col - 2]
}

# Detail rows (2..14).
$data = @(
    @{A=0;  B="000556"; C="国投瑞银新机遇灵活配置混合A";        D="4.50"; E="21.12"; F="1.35"; G="0.0608"; H=2},
    @{A=1;  B="007509"; C="华商润丰灵活配置混合C";              D="2.97"; E="40.93"; F="1.14"; G="0.0339"; H=9},
    @{A=2;  B="519615"; C="银河君尚灵活配置混合I";              D="3.59"; E="35.36"; F="0.80"; G="0.0287"; H=2},
    @{A=3;  B="006429"; C="诺安恒鑫混合";                       D="0.64"; E="70.94"; F="4.07"; G="0.0260"; H=8},
    @{A=4;  B="011243"; C="万家惠裕回报6个月持有期混合A";       D="1.54"; E="27.67"; F="1.34"; G="0.0206"; H=4},
    @{A=5;  B="000557"; C="国投瑞银新机遇灵活配置混合C";        D="1.45"; E="21.12"; F="1.35"; G="0.0196"; H=2},
    @{A=6;  B="519613"; C="银河君尚灵活配置混合A";              D="2.10"; E="35.36"; F="0.80"; G="0.0168"; H=2},
    @{A=7;  B="005053"; C="银河量化价值混合A";                  D="0.10"; E="78.55"; F="1.65"; G="0.0016"; H=10},
    @{A=8;  B="011244"; C="万家惠裕回报6个月持有期混合C";       D="0.12"; E="27.67"; F="1.34"; G="0.0016"; H=4},
    @{A=9;  B="005126"; C="银河量化稳进混合";                   D="0.13"; E="55.69"; F="1.08"; G="0.0014"; H=9},
    @{A=10; B="519614"; C="银河君尚灵活配置混合C";              D="0.17"; E="35.36"; F="0.80"; G="0.0014"; H=2},
    @{A=11; B="003598"; C="华商润丰灵活配置混合A";              D="0.06"; E="40.93"; F="1.14"; G="0.0007"; H=9},
    @{A=12; B="013026"; C="银河量化价值混合C";                  D="0.00"; E="78.55"; F="1.65"; G=0;        H=10}
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rec = $data[$i]

    $q3.Cells.Item($row, 1).Value = $rec.A
    Set-TextCell $q3.Cells.Item($row, 2) $rec.B
    $q3.Cells.Item($row, 3).Value = $rec.C
    Set-TextCell $q3.Cells.Item($row, 4) $rec.D
    Set-TextCell $q3.Cells.Item($row, 5) $rec.E
    Set-TextCell $q3.Cells.Item($row, 6) $rec.F

    if ($rec.G -is [string]) {
        Set-TextCell $q3.Cells.Item($row, 7) $rec.G
    } else {
        $q3.Cells.Item($row, 7).Value = $rec.G
    }

    $q3.Cells.Item($row, 8).Value = $rec.H
}

# Restore "总计" as the active/selected sheet, matching the original file.
$total.Activate()
$total.Range("A1").Select() | Out-Null
